$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45916
$ws.Range("B2").Value = 5880.24471902051
$ws.Range("C2").Value = 5134.64581335464
$ws.Range("D2").Value = 5952
$ws.Range("E2").Value = 6531.07542
$ws.Range("F2").Value = -6.93847856941106
$ws.Range("A3").Value = 45917
$ws.Range("B3").Value = 5880.58079845091
$ws.Range("C3").Value = 5168.11329892286
$ws.Range("D3").Value = 2952
$ws.Range("E3").Value = 6531.204028
$ws.Range("F3").Value = 119.447355352998
$ws.Range("A4").Value = 45918
$ws.Range("B4").Value = 5865.40648573042
$ws.Range("C4").Value = 5157.20031388935
$ws.Range("D4").Value = 2952
$ws.Range("E4").Value = 6512.511375
$ws.Range("F4").Value = 118.846050131622
$ws.Range("A5").Value = 45919
$ws.Range("B5").Value = 5867.49291702243
$ws.Range("C5").Value = 4412.71039212604
$ws.Range("D5").Value = 2952
$ws.Range("E5").Value = 6516.650959
$ws.Range("F5").Value = 87.9111847543169
$ws.Range("A6").Value = 45920
$ws.Range("B6").Value = 1207.12769394495
$ws.Range("C6").Value = 2170.79243289725
$ws.Range("D6").Value = 2952
$ws.Range("E6").Value = 2319.790459
$ws.Range("F6").Value = 13.8106332480122
$ws.Range("A7").Value = 45921
$ws.Range("B7").Value = 1096.63644129135
$ws.Range("C7").Value = 2133.85182172049
$ws.Range("D7").Value = 2952
$ws.Range("E7").Value = 2105.170855
$ws.Range("F7").Value = 7.93275980954763
$ws.Range("A8").Value = 45922
$ws.Range("B8").Value = 5979.99309220113
$ws.Range("C8").Value = 4968.1396944729
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 6710.218759
$ws.Range("F8").Value = 114.43189005299
$ws.Range("A9").Value = 45923
$ws.Range("B9").Value = 5979.99309220113
$ws.Range("C9").Value = 5234.52013447001
$ws.Range("D9").Value = 2952
$ws.Range("E9").Value = 6710.218759
$ws.Range("F9").Value = 125.53107505287
$ws.Range("A10").Value = 45924
$ws.Range("B10").Value = 5979.99309220113
$ws.Range("C10").Value = 5490.02459506182
$ws.Range("D10").Value = 2952
$ws.Range("E10").Value = 6710.218759
$ws.Range("F10").Value = 136.177094244195
$ws.Range("A11").Value = 45925
$ws.Range("B11").Value = 5979.99309220113
$ws.Range("C11").Value = 6276.92091608982
$ws.Range("D11").Value = 2952
$ws.Range("E11").Value = 6710.218759
$ws.Range("F11").Value = 168.964440953696
$ws.Range("A12").Value = 45926
$ws.Range("B12").Value = 5979.99309220113
$ws.Range("C12").Value = 5871.74142202389
$ws.Range("D12").Value = 2952
$ws.Range("E12").Value = 6710.218759
$ws.Range("F12").Value = 152.081962034282
$ws.Range("A13").Value = 45927
$ws.Range("B13").Value = 1198.04914556408
$ws.Range("C13").Value = 3546.91777266219
$ws.Range("D13").Value = 2952
$ws.Range("E13").Value = 2321.483151
$ws.Range("F13").Value = 71.5979907540877
$ws.Range("A14").Value = 45928
$ws.Range("B14").Value = 1070.29087265608
$ws.Range("C14").Value = 3519.63727453148
$ws.Range("D14").Value = 2952
$ws.Range("E14").Value = 2183.373747
$ws.Range("F14").Value = 70.0300062031418
$ws.Range("A15").Value = 45929
$ws.Range("B15").Value = 6235.87706540818
$ws.Range("C15").Value = 6456.79425929188
$ws.Range("D15").Value = 2952
$ws.Range("E15").Value = 7031.626518
$ws.Range("F15").Value = 179.189321328487
